$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.540.31"
$ws.Range("E2").Value = "  -7.70%  "
$ws.Range("D3").Value = "1.686.04"
$ws.Range("E3").Value = "  -6.61%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'216.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.57%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.4956"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -16.65%  "
$ws.Range("D8").Value = "'0.2602"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -6.45%  "
$ws.Range("D9").Value = "'21.59"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -7.99%  "
$ws.Range("D10").Value = "'0.06106"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -10.99%  "
$ws.Range("D11").Value = "'0.07267"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.78%  "
$ws.Range("D12").Value = "1.690.57"
$ws.Range("E12").Value = "  -6.43%  "
$ws.Range("D13").Value = "'4.426"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.49%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.916.20"
$ws.Range("E14").Value = "  -6.56%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.5707"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -9.26%  "
$ws.Range("D16").Value = "'0.000008243"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -11.51%  "
$ws.Range("D17").Value = "'64.48"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -14.61%  "
$ws.Range("D18").Value = "26.608.67"
$ws.Range("E18").Value = "  -7.31%  "
$ws.Range("D19").Value = "'4.996"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -9.01%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'10.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.98%  "
$ws.Range("D22").Value = "'181.93"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -13.92%  "
$ws.Range("D23").Value = "'6.155"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -10.46%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'144.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.49%  "
$ws.Range("D26").Value = "'7.547"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("D27").Value = "'0.1126"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -11.69%  "
$ws.Range("D28").Value = "'15.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.26%  "
$ws.Range("D29").Value = "'1.313"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.10%  "
$ws.Range("D30").Value = "'0.05570"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -10.32%  "
$ws.Range("D31").Value = "'1.323"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.90%  "
$ws.Range("D32").Value = "'3.468"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.44%  "
$ws.Range("D33").Value = "'3.456"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").Value = "'1.642"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("D35").Value = "'1.005"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("D36").Value = "'2.378"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.87%  "
$ws.Range("D37").Value = "'0.5852"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -9.08%  "
$ws.Range("D38").Value = "'2.635"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("D39").Value = "'0.01581"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("D40").Value = "1.069.94"
$ws.Range("E40").Value = "  -6.55%  "
$ws.Range("E41").Value = "  -8.37%  "
$ws.Range("D42").Value = "'0.8493"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'98.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("D45").Value = "1.844.04"
$ws.Range("E45").Value = "  -6.10%  "
$ws.Range("D46").Value = "'56.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.58%  "
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "'8.031"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").Value = "'0.4331"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'0.05213"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.59%  "
